$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to Text format so numeric-looking strings
# (e.g. "66.50", "0.000009937", "29.374.99") are preserved exactly
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.374.99'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.846.23'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '240.17'
$ws.Range('D6').Value = '0.6355'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.07555'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.2965'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').Value = '24.65'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('D11').Value = '0.07734'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '1.845.38'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '4.991'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '83.09'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '0.000009937'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '6.165'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '29.397.60'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').Value = '230.33'
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').Value = '12.45'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '7.554'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '156.93'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = '8.373'
$ws.Range('E26').Value = '  -1.01%  '
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '1.462'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').Value = '0.05705'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').Value = '4.031'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').Value = '1.846'
$ws.Range('E33').Value = '  -2.95%  '
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').Value = '0.7164'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '2.594'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '1.255.28'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('D38').Value = '2.790'
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('D39').Value = '0.01807'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = '0.9079'
$ws.Range('E40').Value = '  -0.65%  '
$ws.Range('D41').Value = '6.199'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = '2.005.51'
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('D44').Value = '101.76'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').Value = '66.50'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('D47').Value = '7.063'
$ws.Range('E47').Value = '  -3.83%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').Value = '0.4024'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.704'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1128'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05737'
$ws.Range('E51').Value = '  -0.21%  '
